$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The table grows by one row (18 data rows -> 19 data rows). Give the new last row (19)
# the same formatting as the current last row (18), i.e. the date number format on
# column A, before shifting any data, so the table's per-row style stays consistent.
$ws.Range("A18").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Shift the existing data rows (2-18) down by one (to 3-19), cell by cell, preserving
# each cell's existing style/format (a plain value copy, no formatting change).
for ($r = 18; $r -ge 2; $r--) {
    $dst = $r + 1
    for ($c = 1; $c -le 5; $c++) {
        $srcCell = $ws.Cells.Item($r, $c)
        $dstCell = $ws.Cells.Item($dst, $c)
        $v = $srcCell.Value2
        if ($null -eq $v) {
            $dstCell.ClearContents()
        } else {
            $dstCell.Value2 = $v
        }
    }
}

# Write the corrected/recomputed data for all rows 2-19:
# columns: A=date_of_forecast, B=y_0, C=y_0_forecast, D=y_1, E=y_1_forecast
$data = @(
    @(39400, 2007, 2.070003986395053, 2008, $null),
    @(39765, 2008, 0.517569958955022, 2009, $null),
    @(40130, 2009, -3.956152295564896, 2010, $null),
    @(40494, 2010, 1.234995474941392, 2011, 0.7327527981400461),
    @(40862, 2011, 0.899360810820804, 2012, 0.4113848771853501),
    @(41228, 2012, 0.9010266119894084, 2013, 1.203907967581297),
    @(41592, 2013, 0.02019328874804938, 2014, -0.03860754389363175),
    @(41957, 2014, 0.1729981757035093, 2015, 0.2629870913912535),
    @(42321, 2015, 0.09752710595589686, 2016, 0.1987429576382649),
    @(42689, 2016, -0.5280591151586633, 2017, 0.1903092973221776),
    @(43053, 2017, 0.07201851318385799, 2018, 0.2843016498274009),
    @(43418, 2018, 0.3727661260635617, 2019, -0.8612142616933327),
    @(43783, 2019, -0.801759526476209, 2020, 0.06491682578968483),
    @(44159, 2020, -1.103489789942047, 2021, 1.323658311025055),
    @(44525, 2021, 0.9704846793491928, 2022, -0.6989646400249128),
    @(44890, 2022, -0.7009264669202708, 2023, -0.0234350458557242),
    @(45254, 2023, 0.3928252664241905, 2024, 0.196134499498668),
    @(45618, 2024, 0.3224026462283813, 2025, -0.6671574593505647)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]

    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]

    if ($null -eq $row[2]) {
        $ws.Cells.Item($r, 3).ClearContents()
    } else {
        $ws.Cells.Item($r, 3).Value2 = $row[2]
    }

    $ws.Cells.Item($r, 4).Value2 = $row[3]

    if ($null -eq $row[4]) {
        $ws.Cells.Item($r, 5).ClearContents()
    } else {
        $ws.Cells.Item($r, 5).Value2 = $row[4]
    }
}
